# Add a new Eurobarometer wave (ZA7902 / EB 97.5, June-July 2022) as the
# most-recent entry in the survey list, i.e. insert a new row 2 (just below
# the header row) and shift all existing data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 2 (the former top data row,
# EB 97.3 / ZA7888), pushing everything else down.
$ws.Rows.Item(2).Insert() | Out-Null

# Fill in the new record. Column order on the sheet is:
#   A = archive_id, B = wave, C = timeframe, D = description
# The wave value is entered with a leading apostrophe so it is stored as
# text (matching how the other "wave" values such as 97.3 / 93.2 are
# stored), not as a number.
$ws.Range("A2").Value = "ZA7902"
$ws.Range("B2").Value = "'97.5"
$ws.Range("D2").Value = "Standard Eurobarometer 97 (COVID-19 Pandemic)"
$ws.Range("C2").Value = "June-July 2022"

# Leave the selection on the timeframe cell that was filled in last.
$ws.Range("C2").Select() | Out-Null
